# Add a new "transfer" row (row 10) to the "Oddziały" sheet, mirroring the
# existing Danielewski Paweł / 2TFB / Techniki fryzjerskie entry (row 6) but
# with new "Przeniesiono z" / "Przeniesiono na" schedule-slot values.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Oddziały")

$ws.Range("A10").Value = "30.01.2026, 9, 14:55-15:40, sala: 5"
$ws.Range("B10").Value = "30.01.2026, 8, 14:05-14:50, sala: 4"
$ws.Range("C10").Value = "Danielewski Paweł"
$ws.Range("D10").Value = "-"
$ws.Range("E10").Value = "2TFB"
$ws.Range("F10").Value = "Techniki fryzjerskie"
$ws.Range("G10").Value = ""
